# Update the October Deskcount worksheet:
#  - Mark several office locations as excluded from the occupancy calculation
#    (column F, "Include in Occupancy Calculation": Yes -> No)
#  - Correct the Melbourne deskcount value (column C)
#  - Leave the sheet scrolled/selected where the user last left off

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Deskcount")

# Rows whose "Include in Occupancy Calculation" flag flips from Yes to No.
$ws.Range("F16").Value = "No"   # Greenwood Village
$ws.Range("F38").Value = "No"   # Tampa
$ws.Range("F43").Value = "No"   # Madrid
$ws.Range("F44").Value = "No"   # Melbourne
$ws.Range("F47").Value = "No"   # Santiago
$ws.Range("F48").Value = "No"   # Sao Paulo

# Corrected deskcount for Melbourne.
$ws.Range("C44").Value = 32

# Restore the view/selection state as last saved by the user.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D45").Select()
